$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; existing rows 27-86 shift down to 28-87
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new record's data
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44581
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100107
$ws.Range("H27").Value = "Otros"
$ws.Range("I27").Value = 100107011
$ws.Range("J27").Value = "Tuna"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 20000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 20000
$ws.Range("Q27").Value = "$/caja 20 kilos"
$ws.Range("R27").Value = "Provincia de Limarí"
$ws.Range("S27").Value = 1000
$ws.Range("T27").Value = 20
